$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain ".Value = <numeric-looking text>" gets auto-converted by Excel
# into a real number (losing the original text formatting/precision).
# For those cells we instead enter a formula that evaluates to the exact
# text, then convert it to a static value via Copy + PasteSpecial(values),
# which keeps the cell as plain text without touching its style/format.

$ws.Range("D2").Value = '64.433.33'
$ws.Range("E2").Value = '  -0.28%  '
$ws.Range("D3").Value = '3.496.36'
$ws.Range("E3").Value = '  -0.31%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Formula = '="586.70"'
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = '  +0.09%  '
$ws.Range("D6").Formula = '="134.89"'
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value = '  +1.96%  '
$ws.Range("D7").Value = '3.496.68'
$ws.Range("E7").Value = '  -0.31%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("E9").Value = '  -0.96%  '
$ws.Range("E10").Value = '  -0.30%  '
$ws.Range("E11").Value = '  -1.72%  '
$ws.Range("D12").Formula = '="0.376"'
$ws.Range("D12").Copy()
$ws.Range("D12").PasteSpecial(-4163)
$ws.Range("E12").Value = '  -3.18%  '
$ws.Range("D13").Value = '4.089.56'
$ws.Range("E13").Value = '  -0.40%  '
$ws.Range("E14").Value = '  +0.57%  '
$ws.Range("E15").Value = '  +1.18%  '
$ws.Range("D16").Value = '3.492.59'
$ws.Range("E16").Value = '  -0.56%  '
$ws.Range("D17").Value = '64.435.70'
$ws.Range("E17").Value = '  -0.31%  '
$ws.Range("D18").Formula = '="25.17"'
$ws.Range("D18").Copy()
$ws.Range("D18").PasteSpecial(-4163)
$ws.Range("E18").Value = '  -9.55%  '
$ws.Range("D19").Formula = '="10.03"'
$ws.Range("D19").Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("E19").Value = '  +0.71%  '
$ws.Range("D20").Formula = '="5.65"'
$ws.Range("D20").Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("E20").Value = '  -0.38%  '
$ws.Range("E21").Value = '  -4.78%  '
$ws.Range("D22").Formula = '="385.96"'
$ws.Range("D22").Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("E22").Value = '  -1.99%  '
$ws.Range("E23").Value = '  -1.81%  '
$ws.Range("D24").Value = '3.635.36'
$ws.Range("E24").Value = '  -0.46%  '
$ws.Range("D25").Formula = '="74.06"'
$ws.Range("D25").Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("E25").Value = '  +0.78%  '
$ws.Range("E26").Value = '  +0.00%  '
$ws.Range("E27").Value = '  +3.47%  '
$ws.Range("E28").Value = '  -0.21%  '
$ws.Range("B29").Value = 'Binance-PegBSC-USD'
$ws.Range("C29").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D29").Formula = '="1.00"'
$ws.Range("D29").Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("E29").Value = '  +0.08%  '
$ws.Range("B30").Value = 'Fetch.AI'
$ws.Range("C30").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D30").Formula = '="1.54"'
$ws.Range("D30").Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("E30").Value = '  -2.85%  '
$ws.Range("E31").Value = '  -1.43%  '
$ws.Range("E32").Value = '  +0.02%  '
$ws.Range("D33").Value = '3.516.21'
$ws.Range("E33").Value = '  +0.14%  '
$ws.Range("E34").Value = '  +0.00%  '
$ws.Range("E35").Value = '  +0.06%  '
$ws.Range("E36").Value = '  -2.34%  '
$ws.Range("D37").Formula = '="5.30"'
$ws.Range("D37").Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("E37").Value = '  +0.96%  '
$ws.Range("E38").Value = '  -3.30%  '
$ws.Range("D39").Formula = '="6.83"'
$ws.Range("D39").Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("E39").Value = '  -2.33%  '
$ws.Range("D40").Formula = '="162.51"'
$ws.Range("D40").Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("E40").Value = '  -5.07%  '
$ws.Range("D41").Formula = '="0.0782"'
$ws.Range("D41").Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("E41").Value = '  -3.09%  '
$ws.Range("D42").Formula = '="0.806"'
$ws.Range("D42").Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("E42").Value = '  -1.20%  '
$ws.Range("D43").Formula = '="25.80"'
$ws.Range("D43").Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("E43").Value = '  -1.74%  '
$ws.Range("D45").Formula = '="41.89"'
$ws.Range("D45").Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("E45").Value = '  -0.51%  '
$ws.Range("B46").Value = 'ONDO'
$ws.Range("C46").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D46").Formula = '="1.21"'
$ws.Range("D46").Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("E46").Value = '  +0.16%  '
$ws.Range("B47").Value = 'Filecoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D47").Formula = '="4.41"'
$ws.Range("D47").Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("E47").Value = '  +0.47%  '
$ws.Range("E48").Value = '  -0.33%  '
$ws.Range("D49").Value = '2.477.37'
$ws.Range("E49").Value = '  +1.33%  '
$ws.Range("E50").Value = '  -1.88%  '
$ws.Range("D51").Formula = '="0.905"'
$ws.Range("D51").Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("E51").Value = '  +1.01%  '

$excel.CutCopyMode = 0
